$wb = $excel.ActiveWorkbook

# Sheet "展览" (Sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1908
$ws1.Range("F5").Value = 183
$ws1.Range("F6").Value = 2767
$ws1.Range("F10").Value = 1581
$ws1.Range("F11").Value = 560
$ws1.Range("F12").Value = 48
$ws1.Range("F13").Value = 342
$ws1.Range("F22").Value = 19
$ws1.Range("F23").Value = 229
$ws1.Range("F25").Value = 1769
$ws1.Range("F28").Value = 88
$ws1.Range("F29").Value = 570
$ws1.Range("F31").Value = 314
$ws1.Range("F32").Value = 456

# Sheet "全部类型" (Sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1908
$ws4.Range("F6").Value = 183
$ws4.Range("F7").Value = 2767
$ws4.Range("F11").Value = 1581
$ws4.Range("F12").Value = 560
$ws4.Range("F13").Value = 48
$ws4.Range("F14").Value = 342
$ws4.Range("F23").Value = 19
$ws4.Range("F24").Value = 229
$ws4.Range("F26").Value = 1769
$ws4.Range("F29").Value = 88
$ws4.Range("F30").Value = 570
$ws4.Range("F32").Value = 314
$ws4.Range("F33").Value = 456
